$d = $word.ActiveDocument

# The document has a single section whose title-page-vs-rest headers/footers
# carry the Pearson/BTEC logo pictures (inline pictures, inserted via
# <wp:inline>). Word's object model exposes the DrawingML "name" (the
# <wp:docPr name="..."/> attribute) only through the *floating* Shape
# object, so each inline picture is temporarily promoted with
# ConvertToShape(), renamed, then demoted back with ConvertToInlineShape()
# so the <wp:inline> wrapper (and its layout) is preserved.

function Rename-InlinePicture($inlineShape, $newName) {
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    $floating.ConvertToInlineShape() | Out-Null
}

$section = $d.Sections.Item(1)

# First-page footer (footer1.xml): PearsonLogo picture, id=3 -> image1.png
$firstFooterShape = $section.Footers.Item(2).Range.InlineShapes.Item(1)
Rename-InlinePicture $firstFooterShape "image1.png"

# Default footer (footer2.xml): PearsonLogo picture, id=2 -> image1.png
$defaultFooterShape = $section.Footers.Item(1).Range.InlineShapes.Item(1)
Rename-InlinePicture $defaultFooterShape "image1.png"

# First-page header (header1.xml): BTec_Logo-Orange picture, id=1 -> image2.jpg
$firstHeaderShape = $section.Headers.Item(2).Range.InlineShapes.Item(1)
Rename-InlinePicture $firstHeaderShape "image2.jpg"
